# VG-FS-ADO-Sync.xlsx: add an "AssignedTo" field mapping.
#
# 1) ProductsFields sheet: append a new field-mapping row for
#    System.AssignedTo -> AssignedTo (FS_TO_ADO direction).
# 2) ProductsData sheet: insert a new "AssignedTo" data column (before the
#    existing WorkItemType column) whose values mirror the Developer column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ProductsFields: add row 8
# ---------------------------------------------------------------------
$wsFields = $wb.Worksheets.Item("ProductsFields")

$wsFields.Range("A8").Value = "System.AssignedTo"
$wsFields.Range("B8").Value = "AssignedTo"
$wsFields.Range("C8").Value = "FS_TO_ADO"

# Match the formatting of the row above (row 7) for the new row.
$wsFields.Range("C7").Copy()
$wsFields.Range("C8").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# ProductsData: insert new column H ("AssignedTo") before WorkItemType
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("ProductsData")

# Shift the existing WorkItemType column (H) one column to the right (I).
$wsData.Columns.Item(8).Insert(-4161)   # xlShiftToRight

# Populate the header for the new column, keeping the header's original
# look (same as the WorkItemType header that was just shifted to I1).
$wsData.Range("I1").Copy()
$wsData.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$wsData.Range("H1").Value = "AssignedTo"

# Populate the new column's data with the same assignee as the Developer
# column (F), copying format along with the values.
$wsData.Range("F2:F11").Copy()
$wsData.Range("H2:H11").PasteSpecial(-4104)   # xlPasteAll

# Restore the explicit custom width for the new column (matches the width
# the old WorkItemType column used to have before it moved to I).
$wsData.Columns.Item(8).ColumnWidth = 21.166666666666668

# ---------------------------------------------------------------------
# Selections: ProductsFields -> C11, ProductsData -> H2:H11 (active tab)
# ---------------------------------------------------------------------
$wsFields.Range("C11").Select()
$wsData.Range("H2:H11").Select()
